# account_bank_statement_import_adyen / adyen_test.xlsx
# [MIG] 12.0 account_bank_statement_import_adyen, account_bank_statement_clearing_account
#
# Content edits reproduced from the canonical-OOXML diff:
#  1. Currency code "EUR" -> "USD" (shared string used throughout columns K/O)
#  2. The "Gross Credit (GC)" value in M10: 666 -> 1598
#  3. The date/time number format is re-cased to lower-case tokens
#     (YYYY-MM-DD HH:MM:SS -> yyyy-mm-dd hh:mm:ss) for the whole date column
#  4. Active selection on the sheet moves to L9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. EUR -> USD, whole-cell match only (exact token, not a substring replace)
$ws.Cells.Replace("EUR", "USD", 1)

# 2. M10 numeric value update
$ws.Range("M10").Value = 1598

# 3. Re-case the custom date/time format applied to the whole "Creation Date"
#    column (G5:G25) so it matches the new lower-case format code.
$ws.Range("G5:G25").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# 4. Update the active selection/active cell to L9
$ws.Range("L9").Select()
